# Auto-generated Excel COM-interop script applying the cryptos.xlsx data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.807.51"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.006.99"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.73%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.60"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.31%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.563"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.011.26"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.112"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.37"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.364"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.529.84"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.68%  "
$ws.Range("E14").Value = "  -3.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.905.88"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.87"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.013.87"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.21%  "
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "395.56"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.09"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.85"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.60"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.96%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.13"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.465"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("E26").Value = "  -5.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0962"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.62"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.44"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "160.93"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.67"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.02"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.20%  "
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("E37").Value = "  -2.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.57"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.484.10"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -9.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.57"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.70%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.90"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.71%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.46"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.663"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0593"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.81%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0246"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.02"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.86"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0944"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "263.18"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.71%  "
